$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts existing rows 16-82 down to 17-83
$ws.Rows("16:16").Insert()

# Populate the new row 16 with the new record's data
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44558
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112030
$ws.Range("G16").Value = "Poroto granado"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 23000
$ws.Range("L16").Value = 23000
$ws.Range("M16").Value = 23000
$ws.Range("N16").Value = "`$/saco 25 kilos"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 920
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
